$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# '27÷5=5, 2' -> '94÷3=31, 1'
$t.Cell(1, 1).Range.Text = "94÷3=31, 1"

# '32÷4=8, 0' -> '47÷5=9, 2'
$t.Cell(1, 2).Range.Text = "47÷5=9, 2"

# '21÷4=5, 1' -> '29÷5=5, 4'
$t.Cell(1, 3).Range.Text = "29÷5=5, 4"

# '23÷3=7, 2' -> '13÷2=6, 1'
$t.Cell(1, 4).Range.Text = "13÷2=6, 1"

# '56÷8=7, 0' -> '27÷2=13, 1'
$t.Cell(1, 5).Range.Text = "27÷2=13, 1"

# '89÷2=44, 1' -> '74÷5=14, 4'
$t.Cell(5, 1).Range.Text = "74÷5=14, 4"

# '23÷2=11, 1' -> '12÷2=6, 0'
$t.Cell(5, 2).Range.Text = "12÷2=6, 0"

# '87÷8=10, 7' -> '37÷2=18, 1'
$t.Cell(5, 3).Range.Text = "37÷2=18, 1"

# '48÷6=8, 0' -> '78÷2=39, 0'
$t.Cell(5, 4).Range.Text = "78÷2=39, 0"

# '52÷2=26, 0' -> '37÷6=6, 1'
$t.Cell(5, 5).Range.Text = "37÷6=6, 1"

# '85÷7=12, 1' -> '21÷6=3, 3'
$t.Cell(9, 1).Range.Text = "21÷6=3, 3"

# '76÷8=9, 4' -> '12÷2=6, 0'
$t.Cell(9, 2).Range.Text = "12÷2=6, 0"

# '69÷4=17, 1' -> '43÷7=6, 1'
$t.Cell(9, 3).Range.Text = "43÷7=6, 1"

# '97÷5=19, 2' -> '41÷3=13, 2'
$t.Cell(9, 4).Range.Text = "41÷3=13, 2"

# '48÷3=16, 0' -> '55÷4=13, 3'
$t.Cell(9, 5).Range.Text = "55÷4=13, 3"

# '70÷4=17, 2' -> '63÷8=7, 7'
$t.Cell(13, 1).Range.Text = "63÷8=7, 7"

# '70÷4=17, 2' -> '71÷6=11, 5'
$t.Cell(13, 2).Range.Text = "71÷6=11, 5"

# '29÷9=3, 2' -> '83÷5=16, 3'
$t.Cell(13, 3).Range.Text = "83÷5=16, 3"

# '45÷6=7, 3' -> '53÷5=10, 3'
$t.Cell(13, 4).Range.Text = "53÷5=10, 3"

# '21÷6=3, 3' -> '40÷8=5, 0'
$t.Cell(13, 5).Range.Text = "40÷8=5, 0"

# '90÷4=22, 2' -> '43÷8=5, 3'
$t.Cell(17, 1).Range.Text = "43÷8=5, 3"

# '46÷6=7, 4' -> '72÷9=8, 0'
$t.Cell(17, 2).Range.Text = "72÷9=8, 0"

# '43÷4=10, 3' -> '58÷4=14, 2'
$t.Cell(17, 3).Range.Text = "58÷4=14, 2"

# '38÷3=12, 2' -> '14÷5=2, 4'
$t.Cell(17, 4).Range.Text = "14÷5=2, 4"

# '41÷7=5, 6' -> '39÷5=7, 4'
$t.Cell(17, 5).Range.Text = "39÷5=7, 4"
